$wb = $excel.ActiveWorkbook

# --- Sheet "계획표" (plan sheet, 1st sheet) ---
$ws1 = $wb.Worksheets.Item(1)

# Row 7: fill in the E7 date (240524) - centered like the rest of the date columns
$e7 = $ws1.Cells.Item(7, 5)
$e7.Value2 = 240524
$e7.HorizontalAlignment = -4108   # xlCenter
$e7.VerticalAlignment = -4108     # xlCenter

# Row 10: update the problem text in B10 and add the completion date in C10
$b10 = $ws1.Cells.Item(10, 2)
$b10.Value2 = "1260 (DFS/BFS 뼈대문제임 잘 숙지하기!)"

$c10 = $ws1.Cells.Item(10, 3)
$c10.Value2 = 240524
$c10.HorizontalAlignment = -4108  # xlCenter
$c10.VerticalAlignment = -4108    # xlCenter

# Update the active selection to reflect where the editing ended up
$ws1.Range("E8").Select() | Out-Null

$wb.Save()
